$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at row 17 (this pushes the existing rows 17-37 down to 18-38,
# and Excel automatically updates the formulas/ranges that reference them)
$ws.Rows.Item(17).Insert()

# Fill in the new row 17 with the "REST API für jedes Table" task
$ws.Range("A17").Value = "2.4"
$ws.Range("B17").Value = "REST API für jedes Table"
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 2.3
$ws.Range("E17").Value = "NK"

# Renumber / update the rows that used to be 17-24 and are now 18-25,
# reflecting the new dependency numbers and adjusted effort (hours)
$ws.Range("A18").Value = "2.5"
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 2.4

$ws.Range("A19").Value = "2.6"
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 2.4

$ws.Range("A20").Value = "2.7"
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 2.5

$ws.Range("A21").Value = "2.8"
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 2.4

$ws.Range("A22").Value = "2.9"
$ws.Range("C22").Value = 8

$ws.Range("A23").Value = "2.10"
$ws.Range("B23").Value = "Zimmersuche GUI"

$ws.Range("A24").Value = "2.11"
$ws.Range("B24").Value = "Gastdetails GUI"
$ws.Range("D24").Value = 2.7

$ws.Range("A25").Value = "2.12"

# Set the active cell / selection to match the author's final view
$ws.Range("A17:E17").Select()
